$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume/change (E) columns with latest scraped values.
# Numeric-looking price strings are written with a temporary Text format so Excel
# keeps them as literal strings (preserving trailing zeros / exact formatting),
# then the cell style is reset back to Normal so no stray formatting is introduced.

$ws.Range("D2").Value = "60.024.40"
$ws.Range("E2").Value = "  +3.39%  "
$ws.Range("D3").Value = "2.419.12"
$ws.Range("E3").Value = "  +2.79%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "551.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.81%  "
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.581"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.25%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.54%  "
$ws.Range("D14").Value = "2.851.28"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("D15").Value = "59.983.95"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "2.422.30"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "331.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.08%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.58%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.87%  "
$ws.Range("E28").Value = "  +6.01%  "
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.02"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +4.80%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.415"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "314.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.578"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.407"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.64"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "

Write-Host "Applied crypto list updates"
